$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (Q1:V1) do not exist yet; give them the same
# bold/border/center-top header style as the existing header row before
# filling in their text ---
$ws.Range("P1").Copy()
$ws.Range("Q1:V1").PasteSpecial(-4122)

# --- Update header row (1): rename existing + add new headers ---
$ws.Range("H1").Value = "arrecadado_avg"
$ws.Range("I1").Value = "arrecadado_std"
$ws.Range("J1").Value = "arrecadado_min"
$ws.Range("K1").Value = "arrecadado_max"
$ws.Range("M1").Value = "apoio_std"
$ws.Range("N1").Value = "apoio_min"
$ws.Range("O1").Value = "apoio_max"
$ws.Range("P1").Value = "contribuicoes"
$ws.Range("Q1").Value = "contribuicoes_med"
$ws.Range("R1").Value = "contribuicoes_std"
$ws.Range("S1").Value = "contribuicoes_min"
$ws.Range("T1").Value = "contribuicoes_max"
$ws.Range("U1").Value = "menor_ano"
$ws.Range("V1").Value = "maior_ano"

# --- New columns get number format matching their group (R$ for apoio_*, #,##0 for contribuicoes_*) ---
$ws.Range("M2:O22").NumberFormat = "R$ #,##0.00"
$ws.Range("P2:Q22").NumberFormat = "#,##0"
$ws.Range("R2:T22").NumberFormat = "#,##0"

# --- Update data rows (2-22) ---
# Row 2: angelo_agostini
$ws.Range("L2").Value = 78.57672143399174
$ws.Range("M2").Value = 23.63877229505447
$ws.Range("N2").Value = 32.2694534583262
$ws.Range("O2").Value = 151.1292159501072
$ws.Range("P2").Value = 35293
$ws.Range("Q2").Value = 578.5737704918033
$ws.Range("R2").Value = 929.9762624034444
$ws.Range("S2").Value = 55
$ws.Range("T2").Value = 6494
$ws.Range("U2").Value = 2013
$ws.Range("V2").Value = 2023

# Row 3: ccxp
$ws.Range("L3").Value = 87.11000777143747
$ws.Range("M3").Value = 35.07795799700576
$ws.Range("N3").Value = 33.80063482849972
$ws.Range("O3").Value = 257.7853211115706
$ws.Range("P3").Value = 34535
$ws.Range("Q3").Value = 274.0873015873016
$ws.Range("R3").Value = 239.4804716828918
$ws.Range("S3").Value = 36
$ws.Range("T3").Value = 1815
$ws.Range("U3").Value = 2014
$ws.Range("V3").Value = 2023

# Row 4: disputa
$ws.Range("L4").Value = 95.92481440598637
$ws.Range("M4").Value = 62.28205825061158
$ws.Range("N4").Value = 21.61624650544615
$ws.Range("O4").Value = 792.0360759681182
$ws.Range("P4").Value = 70527
$ws.Range("Q4").Value = 287.865306122449
$ws.Range("R4").Value = 263.6097661931087
$ws.Range("S4").Value = 11
$ws.Range("T4").Value = 1588
$ws.Range("U4").Value = 2012
$ws.Range("V4").Value = 2023

# Row 5: erotismo
$ws.Range("L5").Value = 94.90136016654937
$ws.Range("M5").Value = 42.06922775736284
$ws.Range("N5").Value = 37.30506273801686
$ws.Range("O5").Value = 323.2845357010965
$ws.Range("P5").Value = 25810
$ws.Range("Q5").Value = 314.7560975609756
$ws.Range("R5").Value = 268.2786261681259
$ws.Range("S5").Value = 22
$ws.Range("T5").Value = 1539
$ws.Range("U5").Value = 2012
$ws.Range("V5").Value = 2023

# Row 6: fantasia
$ws.Range("L6").Value = 88.02163529111277
$ws.Range("M6").Value = 33.79527781109714
$ws.Range("N6").Value = 23.15006403629383
$ws.Range("O6").Value = 213.9734252018395
$ws.Range("P6").Value = 54018
$ws.Range("Q6").Value = 300.1
$ws.Range("R6").Value = 283.7555817382237
$ws.Range("S6").Value = 3
$ws.Range("T6").Value = 1711
$ws.Range("U6").Value = 2012
$ws.Range("V6").Value = 2023

# Row 7: ficcao_cientifica
$ws.Range("L7").Value = 94.68234986858575
$ws.Range("M7").Value = 69.67419713686293
$ws.Range("N7").Value = 30.56837093393595
$ws.Range("O7").Value = 792.0360759681182
$ws.Range("P7").Value = 61529
$ws.Range("Q7").Value = 343.7374301675978
$ws.Range("R7").Value = 493.6165473708339
$ws.Range("S7").Value = 1
$ws.Range("T7").Value = 5879
$ws.Range("U7").Value = 2012
$ws.Range("V7").Value = 2023

# Row 8: fiq
$ws.Range("L8").Value = 86.76189291082824
$ws.Range("M8").Value = 32.20538078565799
$ws.Range("N8").Value = 21.61624650544615
$ws.Range("O8").Value = 199.8601709743299
$ws.Range("P8").Value = 54948
$ws.Range("Q8").Value = 339.1851851851852
$ws.Range("R8").Value = 351.9795170180732
$ws.Range("S8").Value = 31
$ws.Range("T8").Value = 3266
$ws.Range("U8").Value = 2011
$ws.Range("V8").Value = 2023

# Row 9: folclore
$ws.Range("L9").Value = 90.34545930449829
$ws.Range("M9").Value = 67.62639234759273
$ws.Range("N9").Value = 32.2694534583262
$ws.Range("O9").Value = 792.0360759681182
$ws.Range("P9").Value = 51967
$ws.Range("Q9").Value = 371.1928571428571
$ws.Range("R9").Value = 416.4698145405411
$ws.Range("S9").Value = 17
$ws.Range("T9").Value = 3266
$ws.Range("U9").Value = 2012
$ws.Range("V9").Value = 2023

# Row 10: herois
$ws.Range("L10").Value = 96.01288922432792
$ws.Range("M10").Value = 67.13300946444733
$ws.Range("N10").Value = 21.61624650544615
$ws.Range("O10").Value = 792.0360759681182
$ws.Range("P10").Value = 38102
$ws.Range("Q10").Value = 242.687898089172
$ws.Range("R10").Value = 251.542596013111
$ws.Range("S10").Value = 11
$ws.Range("T10").Value = 1588
$ws.Range("U10").Value = 2012
$ws.Range("V10").Value = 2023

# Row 11: hqmix
$ws.Range("L11").Value = 81.01947145980097
$ws.Range("M11").Value = 26.64100216145226
$ws.Range("N11").Value = 21.61624650544615
$ws.Range("O11").Value = 172.1642729447236
$ws.Range("P11").Value = 42155
$ws.Range("Q11").Value = 390.3240740740741
$ws.Range("R11").Value = 706.124751741857
$ws.Range("S11").Value = 12
$ws.Range("T11").Value = 6494
$ws.Range("U11").Value = 2013
$ws.Range("V11").Value = 2023

# Row 12: hqmix
$ws.Range("L12").Value = 81.01947145980097
$ws.Range("M12").Value = 26.64100216145226
$ws.Range("N12").Value = 21.61624650544615
$ws.Range("O12").Value = 172.1642729447236
$ws.Range("P12").Value = 42155
$ws.Range("Q12").Value = 390.3240740740741
$ws.Range("R12").Value = 706.124751741857
$ws.Range("S12").Value = 12
$ws.Range("T12").Value = 6494
$ws.Range("U12").Value = 2013
$ws.Range("V12").Value = 2023

# Row 13: jogos
$ws.Range("L13").Value = 91.15737775390478
$ws.Range("M13").Value = 33.38796883158906
$ws.Range("N13").Value = 35.29658989882071
$ws.Range("O13").Value = 234.8710142410997
$ws.Range("P13").Value = 61093
$ws.Range("Q13").Value = 303.9452736318408
$ws.Range("R13").Value = 494.1239844310836
$ws.Range("S13").Value = 26
$ws.Range("T13").Value = 6494
$ws.Range("U13").Value = 2012
$ws.Range("V13").Value = 2023

# Row 14: lgbtqiamais
$ws.Range("L14").Value = 88.86048522662777
$ws.Range("M14").Value = 39.11867998547866
$ws.Range("N14").Value = 37.30506273801686
$ws.Range("O14").Value = 245.6155654729304
$ws.Range("P14").Value = 17873
$ws.Range("Q14").Value = 308.1551724137931
$ws.Range("R14").Value = 327.6476927842036
$ws.Range("S14").Value = 8
$ws.Range("T14").Value = 1539
$ws.Range("U14").Value = 2013
$ws.Range("V14").Value = 2023

# Row 15: midia_independente
$ws.Range("L15").Value = 103.6833090125289
$ws.Range("M15").Value = 51.99423496130147
$ws.Range("N15").Value = 36.80839302979295
$ws.Range("O15").Value = 323.2845357010965
$ws.Range("P15").Value = 38250
$ws.Range("Q15").Value = 394.3298969072165
$ws.Range("R15").Value = 391.0591046986345
$ws.Range("S15").Value = 29
$ws.Range("T15").Value = 1711
$ws.Range("U15").Value = 2012
$ws.Range("V15").Value = 2023

# Row 16: politica
$ws.Range("L16").Value = 95.62742936226397
$ws.Range("M16").Value = 44.41979710275996
$ws.Range("N16").Value = 33.80063482849972
$ws.Range("O16").Value = 362.0414364166904
$ws.Range("P16").Value = 51893
$ws.Range("Q16").Value = 336.9675324675325
$ws.Range("R16").Value = 298.7578952675915
$ws.Range("S16").Value = 1
$ws.Range("T16").Value = 1588
$ws.Range("U16").Value = 2012
$ws.Range("V16").Value = 2023

# Row 17: questoes_genero
$ws.Range("L17").Value = 88.50596984923921
$ws.Range("M17").Value = 44.08134994103607
$ws.Range("N17").Value = 46.83761258476419
$ws.Range("O17").Value = 245.6155654729304
$ws.Range("P17").Value = 7569
$ws.Range("Q17").Value = 315.375
$ws.Range("R17").Value = 313.5216257016551
$ws.Range("S17").Value = 51
$ws.Range("T17").Value = 1489
$ws.Range("U17").Value = 2013
$ws.Range("V17").Value = 2023

# Row 18: religiosidade
$ws.Range("L18").Value = 87.04136070513805
$ws.Range("M18").Value = 35.12732957218113
$ws.Range("N18").Value = 21.61624650544615
$ws.Range("O18").Value = 301.8203940790075
$ws.Range("P18").Value = 75649
$ws.Range("Q18").Value = 358.5260663507109
$ws.Range("R18").Value = 555.7348485891068
$ws.Range("S18").Value = 4
$ws.Range("T18").Value = 6494
$ws.Range("U18").Value = 2012
$ws.Range("V18").Value = 2023

# Row 19: saloes_humor
$ws.Range("L19").Value = 75.00500439706998
$ws.Range("M19").Value = 43.45490371809562
$ws.Range("N19").Value = 32.2694534583262
$ws.Range("O19").Value = 185.9579322823807
$ws.Range("P19").Value = 2243
$ws.Range("Q19").Value = 186.9166666666667
$ws.Range("R19").Value = 124.1256389406485
$ws.Range("S19").Value = 11
$ws.Range("T19").Value = 411
$ws.Range("U19").Value = 2013
$ws.Range("V19").Value = 2023

# Row 20: terror
$ws.Range("L20").Value = 91.19226916374137
$ws.Range("M20").Value = 57.41534544136429
$ws.Range("N20").Value = 21.61624650544615
$ws.Range("O20").Value = 792.0360759681182
$ws.Range("P20").Value = 108579
$ws.Range("Q20").Value = 332.045871559633
$ws.Range("R20").Value = 527.8918493761878
$ws.Range("S20").Value = 16
$ws.Range("T20").Value = 6494
$ws.Range("U20").Value = 2012
$ws.Range("V20").Value = 2023

# Row 21: webformatos
$ws.Range("L21").Value = 83.093358381024
$ws.Range("M21").Value = 29.94469803232121
$ws.Range("N21").Value = 21.61624650544615
$ws.Range("O21").Value = 194.2230576381307
$ws.Range("P21").Value = 25369
$ws.Range("Q21").Value = 285.0449438202247
$ws.Range("R21").Value = 251.1921534544801
$ws.Range("S21").Value = 38
$ws.Range("T21").Value = 1547
$ws.Range("U21").Value = 2012
$ws.Range("V21").Value = 2023

# Row 22: zine
$ws.Range("L22").Value = 82.11629402200307
$ws.Range("M22").Value = 32.47932839143289
$ws.Range("N22").Value = 32.2694534583262
$ws.Range("O22").Value = 245.6155654729304
$ws.Range("P22").Value = 28547
$ws.Range("Q22").Value = 254.8839285714286
$ws.Range("R22").Value = 226.6887672930878
$ws.Range("S22").Value = 1
$ws.Range("T22").Value = 1588
$ws.Range("U22").Value = 2012
$ws.Range("V22").Value = 2023
